# Applies updated cryptocurrency price/volume data to Sheet1
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.436.92"
$ws.Range("E2").Value = "  +1.09%  "
# Row 3
$ws.Range("D3").Value = "3.449.78"
$ws.Range("E3").Value = "  +1.98%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.82"
$ws.Range("D5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.45%  "
# Row 7
$ws.Range("D7").Value = "3.451.38"
$ws.Range("E7").Value = "  +2.12%  "
# Row 8
$ws.Range("E8").Value = "  +0.06%  "
# Row 9
$ws.Range("E9").Value = "  +0.67%  "
# Row 10
$ws.Range("E10").Value = "  +3.12%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.10%  "
# Row 13
$ws.Range("D13").Value = "4.042.75"
$ws.Range("E13").Value = "  +2.11%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.34%  "
# Row 15
$ws.Range("E15").Value = "  -0.50%  "
# Row 16
$ws.Range("E16").Value = "  +1.40%  "
# Row 17
$ws.Range("D17").Value = "3.449.69"
$ws.Range("E17").Value = "  +2.08%  "
# Row 18
$ws.Range("D18").Value = "61.582.58"
$ws.Range("E18").Value = "  +1.01%  "
# Row 19
$ws.Range("E19").Value = "  +8.65%  "
# Row 20
$ws.Range("E20").Value = "  +2.44%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.54%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.569"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "
# Row 24
$ws.Range("D24").Value = "3.595.12"
$ws.Range("E24").Value = "  +2.42%  "
# Row 25
$ws.Range("E25").Value = "  +2.24%  "
# Row 26
$ws.Range("E26").Value = "  -0.10%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "
# Row 28
$ws.Range("E28").Value = "  -1.82%  "
# Row 29
$ws.Range("E29").Value = "  +7.31%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.98%  "
# Row 31
$ws.Range("E31").Value = "  +0.01%  "
# Row 32
$ws.Range("E32").Value = "  -13.63%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.44%  "
# Row 34
$ws.Range("E34").Value = "  +0.95%  "
# Row 35
$ws.Range("E35").Value = "  +0.03%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.98%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.76%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "
# Row 39
$ws.Range("E39").Value = "  +2.46%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0790"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.55%  "
# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.53%  "
# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.797"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.19%  "
# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.43%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.85%  "
# Row 47
$ws.Range("E47").Value = "  +1.41%  "
# Row 48
$ws.Range("D48").Value = "2.604.19"
$ws.Range("E48").Value = "  +9.30%  "
# Row 49
$ws.Range("E49").Value = "  -3.82%  "
# Row 50
$ws.Range("E50").Value = "  +2.98%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
